# chore: simulator full-month coverage, persist logs, fix employees
#
# The timekeeping simulator regenerated this employee's export:
#  - several client names on the "Weekly Timesheet" sheet were replaced
#  - hours/rate/total were recalculated (rate is now $88/hr, hours bumped to 9
#    on most days, and 01-16 split into a Regular entry + a new OT entry)
#  - the subtotal / hourly / admin / grand-total rows shifted down to make
#    room for the new OT row, and totals were recalculated
#  - the mirrored flat log ("Jason Schema") got the same new OT row appended
#  - the employee's internal ID was regenerated

$wb = $excel.ActiveWorkbook
$tsSheet  = $wb.Worksheets.Item("Weekly Timesheet")
$logSheet = $wb.Worksheets.Item("Jason Schema")

# ---------------------------------------------------------------------------
# 1) Weekly Timesheet: update the five existing daily rows (client, hours,
#    rate, total). Rate is now 88/hr across the board.
# ---------------------------------------------------------------------------

$tsSheet.Range("B2").Value2 = "Prezzano"
$tsSheet.Range("C2").Value2 = 9
$tsSheet.Range("E2").Value2 = 88
$tsSheet.Range("F2").Value2 = 792

$tsSheet.Range("B3").Value2 = "Vincent"
$tsSheet.Range("C3").Value2 = 9
$tsSheet.Range("E3").Value2 = 88
$tsSheet.Range("F3").Value2 = 792

$tsSheet.Range("B4").Value2 = "Zygmunt"
$tsSheet.Range("C4").Value2 = 9
$tsSheet.Range("E4").Value2 = 88
$tsSheet.Range("F4").Value2 = 792

$tsSheet.Range("B5").Value2 = "Ricca"
$tsSheet.Range("C5").Value2 = 9
$tsSheet.Range("E5").Value2 = 88
$tsSheet.Range("F5").Value2 = 792

$tsSheet.Range("B6").Value2 = "Varricchio"
$tsSheet.Range("C6").Value2 = 4
$tsSheet.Range("E6").Value2 = 88
$tsSheet.Range("F6").Value2 = 352

# ---------------------------------------------------------------------------
# 2) Insert a new row for the 2026-01-16 OT entry right after the existing
#    2026-01-16 Regular row. Inserting here pushes the blank separator row
#    and the SUBTOTAL/Category/HOURLY/ADMIN/GRAND rows down by one, which is
#    exactly what the new layout needs.
# ---------------------------------------------------------------------------

$tsSheet.Rows.Item(7).Insert()

# Copy the date + client text over from the row above so it stays a shared
# text value (2026-01-16 / Varricchio) instead of Excel re-interpreting the
# literal string as a date serial number.
$tsSheet.Range("A6:B6").Copy()
$tsSheet.Range("A7").PasteSpecial(-4163)

$tsSheet.Range("C7").Value2 = 5
$tsSheet.Range("D7").Value2 = "OT"
$tsSheet.Range("E7").Value2 = 88
$tsSheet.Range("F7").Value2 = 660

# ---------------------------------------------------------------------------
# 3) Recalculate the subtotal block (now at rows 9/10/12/13/14 after the
#    insert): total hours, the Reg/OT summary text, and the total dollars.
# ---------------------------------------------------------------------------

$tsSheet.Range("C9").Value2 = 45
$tsSheet.Range("D9").Value2 = "Reg: 40 / OT: 5"
$tsSheet.Range("F9").Value2 = 4180

$tsSheet.Range("F12").Value2 = 4180
$tsSheet.Range("F14").Value2 = 4180

# ---------------------------------------------------------------------------
# 4) Jason Schema (flat log): mirror the same five row updates, then append
#    the new OT row so every timesheet entry is persisted to the log too.
# ---------------------------------------------------------------------------

$logSheet.Range("D2").Value2 = "Prezzano"
$logSheet.Range("E2").Value2 = 9
$logSheet.Range("F2").Value2 = 88
$logSheet.Range("G2").Value2 = 792

$logSheet.Range("D3").Value2 = "Vincent"
$logSheet.Range("E3").Value2 = 9
$logSheet.Range("F3").Value2 = 88
$logSheet.Range("G3").Value2 = 792

$logSheet.Range("D4").Value2 = "Zygmunt"
$logSheet.Range("E4").Value2 = 9
$logSheet.Range("F4").Value2 = 88
$logSheet.Range("G4").Value2 = 792

$logSheet.Range("D5").Value2 = "Ricca"
$logSheet.Range("E5").Value2 = 9
$logSheet.Range("F5").Value2 = 88
$logSheet.Range("G5").Value2 = 792

$logSheet.Range("D6").Value2 = "Varricchio"
$logSheet.Range("E6").Value2 = 4
$logSheet.Range("F6").Value2 = 88
$logSheet.Range("G6").Value2 = 352

# Append the new OT log row by copying row 6's layout/format, then fixing up
# the values that differ (Hours, Total, Type).
$logSheet.Range("A6:I6").Copy()
$logSheet.Range("A7").PasteSpecial(-4163)

$logSheet.Range("E7").Value2 = 5
$logSheet.Range("G7").Value2 = 660
$logSheet.Range("H7").Value2 = "OT"

# ---------------------------------------------------------------------------
# 5) Regenerate the employee's internal ID (every row on the log sheet
#    references the same employee, so update them all).
# ---------------------------------------------------------------------------

$logSheet.Range("B2").Value2 = "emp_emnnysju"
$logSheet.Range("B3").Value2 = "emp_emnnysju"
$logSheet.Range("B4").Value2 = "emp_emnnysju"
$logSheet.Range("B5").Value2 = "emp_emnnysju"
$logSheet.Range("B6").Value2 = "emp_emnnysju"
$logSheet.Range("B7").Value2 = "emp_emnnysju"

Write-Host "edit complete"
